$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = 45987
$ws.Range("A82").Style = $ws.Range("A81").Style
$ws.Range("A82").NumberFormat = $ws.Range("A81").NumberFormat

$ws.Range("B82").Value = "15,3635"
$ws.Range("C82").Value = "15,8798"
$ws.Range("D82").Value = "15,3635"
$ws.Range("E82").Value = "15,3635"
